$p = $ppt.ActivePresentation

# Slide 1 title: merge "First" + " " + "slide" runs into a single run "First slide".
# Setting TextRange.Text directly to the same concatenated value is a no-op for the
# underlying run structure, so first set it to a placeholder to force a rewrite.
$s1 = $p.Slides.Item(1)
$t1 = $s1.Shapes.Item(1).TextFrame.TextRange
$t1.Text = "~"
$t1.Text = "First slide"

# Slide 3 title: merge "Third" + " " + "slide" runs into a single run "Third slide".
$s3 = $p.Slides.Item(3)
$t3 = $s3.Shapes.Item(1).TextFrame.TextRange
$t3.Text = "~"
$t3.Text = "Third slide"
